# Add a new "Greece" market sheet, copied from the existing "Croatia" sheet
# (same layout/styles), then filled in with Greece-specific test data.
$wb = $excel.ActiveWorkbook
$croatia = $wb.Worksheets.Item("Croatia")

# Copy "Croatia" and place the copy right after it (becomes the new last sheet).
$croatia.Copy($null, $croatia)
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "Greece"

# Greece-specific values (market name + Jira reference).
$new.Range("B2").Value2 = "Greece Market"
$new.Range("B4").Value2 = "NGC-4119/T3165"

# Restore the source sheet's selection to "select all" and make the new
# "Greece" sheet the active/selected tab.
$croatia.Activate()
$croatia.Cells.Select()
$new.Activate()
